$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently stored in column L (rows 15-18) before we
# clear them, so we can re-create them further down the sheet.
$wildcardConfigValue = $ws.Range("L15").Value2      # "Data GlobalConfiguration config"
$wildcardValue       = $ws.Range("L16").Value2      # "wildcard"

# Copy formats (fill/border/alignment/font) from the existing styled rows
# so the new rows 28/29 keep a consistent look with rows 16/17.
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null

$ws.Range("C17").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Remove the old column L values - they are being relocated.
$ws.Range("L15").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("L18").ClearContents()

# Re-create the data further down the sheet, in column C.
$ws.Range("C27").Value = $wildcardConfigValue
$ws.Range("C28").Value = $wildcardValue
$ws.Range("C29").Value = "Wildcard"
$ws.Range("C30").Value = $false

# Update the view state to match the new layout.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C27").Select()
